$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.986.14'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.858.01'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.56'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5090'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +1.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3814'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08267'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -7.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.110'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.49'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.199'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -2.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.54'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.858.20'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.196'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001098'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.56'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06602'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.68'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.018'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.000.65'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -4.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.241'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.554'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +1.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.072.05'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.82'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("E29").Value = '  -1.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.30'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1054'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.619'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.600'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.452'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06515'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02411'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2168'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.205'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6449'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +1.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.234'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -4.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.873'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -0.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.16'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -3.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6082'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.09'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.276'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.651'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.996'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.01'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.76'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.05%  '
